$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (preserving the original "inline string" /
# text cell type) even when the value looks like a number (e.g. "0.9998").
# Excel's normal Value-assignment auto-converts number-looking strings to
# real numbers, so we briefly force a Text number format while assigning,
# then restore the cell to the default "Normal" style so no stray
# formatting is left behind (matching the source workbook, whose cells
# carry no explicit style).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "25.856.54"
$ws.Range("E2").Value = "  +0.18%  "
Set-TextValue "D3" "1.743.13"
Set-TextValue "D4" "0.9998"
Set-TextValue "D5" "225.36"
$ws.Range("E5").Value = "  -5.04%  "
Set-TextValue "D6" "0.9998"
$ws.Range("E6").Value = "  -0.07%  "
Set-TextValue "D7" "0.5146"
$ws.Range("E7").Value = "  +1.50%  "
Set-TextValue "D8" "0.2786"
$ws.Range("E8").Value = "  +4.88%  "
Set-TextValue "D9" "38.96"
$ws.Range("E9").Value = "  -5.55%  "
Set-TextValue "D10" "0.06092"
$ws.Range("E10").Value = "  -1.90%  "
Set-TextValue "D11" "1.736.82"
$ws.Range("E11").Value = "  -0.99%  "
Set-TextValue "D12" "0.06990"
$ws.Range("E12").Value = "  +0.70%  "
Set-TextValue "D13" "15.25"
$ws.Range("E13").Value = "  -2.18%  "
Set-TextValue "D14" "0.6358"
$ws.Range("E14").Value = "  +4.74%  "
Set-TextValue "D15" "4.510"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("E16").Value = "  -1.36%  "
Set-TextValue "D17" "0.9995"
$ws.Range("E17").Value = "  -0.16%  "
Set-TextValue "D18" "0.9993"
$ws.Range("E18").Value = "  -0.12%  "
Set-TextValue "D19" "25.875.48"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  -1.53%  "
Set-TextValue "D21" "0.000006591"
$ws.Range("E21").Value = "  -3.57%  "
Set-TextValue "D22" "1.957.34"
$ws.Range("E22").Value = "  -0.97%  "
Set-TextValue "D23" "4.091"
$ws.Range("E23").Value = "  +0.52%  "
Set-TextValue "D24" "8.516"
$ws.Range("E24").Value = "  +4.30%  "
Set-TextValue "D25" "5.106"
$ws.Range("E25").Value = "  -1.68%  "
Set-TextValue "D26" "137.47"
$ws.Range("E26").Value = "  -0.29%  "
Set-TextValue "D27" "1.501"
$ws.Range("E27").Value = "  +2.98%  "
Set-TextValue "D28" "1.820"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  -0.18%  "
Set-TextValue "D30" "102.67"
$ws.Range("E30").Value = "  -0.10%  "
Set-TextValue "D31" "0.08273"
$ws.Range("E31").Value = "  +0.53%  "
Set-TextValue "D32" "3.626"
$ws.Range("E32").Value = "  -1.70%  "
Set-TextValue "D33" "3.392"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  +0.90%  "
Set-TextValue "D35" "2.617"
$ws.Range("E35").Value = "  -1.45%  "
Set-TextValue "D36" "0.9702"
$ws.Range("E36").Value = "  -3.23%  "
Set-TextValue "D37" "0.5983"
$ws.Range("E37").Value = "  -1.37%  "
Set-TextValue "D38" "2.671"
$ws.Range("E38").Value = "  -2.08%  "
Set-TextValue "D39" "0.01554"
$ws.Range("E39").Value = "  +0.21%  "
Set-TextValue "D40" "1.916"
$ws.Range("E40").Value = "  -1.10%  "
Set-TextValue "D41" "0.9990"
$ws.Range("E41").Value = "  -0.18%  "
Set-TextValue "D42" "101.02"
$ws.Range("E42").Value = "  -2.13%  "
Set-TextValue "D43" "0.3829"
$ws.Range("E43").Value = "  -0.13%  "
Set-TextValue "D44" "0.7275"
$ws.Range("E44").Value = "  -1.56%  "
Set-TextValue "D45" "4.871"
$ws.Range("E45").Value = "  -0.48%  "
Set-TextValue "D46" "0.05465"
$ws.Range("E46").Value = "  -0.50%  "
Set-TextValue "D47" "6.242"
$ws.Range("E47").Value = "  +4.82%  "
Set-TextValue "D48" "0.1103"
$ws.Range("E48").Value = "  +2.07%  "
Set-TextValue "D49" "52.21"
Set-TextValue "D50" "29.69"
$ws.Range("E50").Value = "  -1.21%  "
Set-TextValue "D51" "7.529"
$ws.Range("E51").Value = "  -0.73%  "
